$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44874
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = 1050
$ws.Range("P2").Value = 1050

# Row 3
$ws.Range("D3").Value = 44832
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = 2100
$ws.Range("P3").Value = 2100

# Row 4
$ws.Range("D4").Value = 44511
$ws.Range("J4").Value = 500

# Row 5
$ws.Range("D5").Value = 44508
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 950
$ws.Range("P5").Value = 950

# Row 6
$ws.Range("D6").Value = 44845
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 1300
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1400
$ws.Range("P6").Value = 1400

# Row 7
$ws.Range("D7").Value = 44897
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 950
$ws.Range("P7").Value = 950

# Row 8
$ws.Range("D8").Value = 44512
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 900
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 950
$ws.Range("P8").Value = 950

# Row 10
$ws.Range("D10").Value = 44460
$ws.Range("H10").Value = "Verde"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 2200
$ws.Range("L10").Value = 2300
$ws.Range("M10").Value = 2250
$ws.Range("P10").Value = 2250

# Row 11
$ws.Range("D11").Value = 44837
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1900
$ws.Range("P11").Value = 1900

# Row 12
$ws.Range("D12").Value = 44882
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1100
$ws.Range("M12").Value = 1050
$ws.Range("O12").Value = "Región de Ñuble"
$ws.Range("P12").Value = 1050

# Row 13
$ws.Range("D13").Value = 44503
$ws.Range("J13").Value = 400

# Row 14
$ws.Range("D14").Value = 44893
$ws.Range("J14").Value = 2000
$ws.Range("O14").Value = "Región de Ñuble"

# Row 15
$ws.Range("D15").Value = 44855
$ws.Range("J15").Value = 800
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1100
$ws.Range("P15").Value = 1100

# Row 16
$ws.Range("D16").Value = 44516
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("P16").Value = 950

# Row 18
$ws.Range("D18").Value = 44532
$ws.Range("J18").Value = 240
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 900
$ws.Range("M18").Value = 850
$ws.Range("O18").Value = "Provincia de Diguillín"
$ws.Range("P18").Value = 850

# Row 19
$ws.Range("D19").Value = 44523

# Row 20
$ws.Range("D20").Value = 44847
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 1300
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = 1400
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 1400

# Row 21
$ws.Range("D21").Value = 44876
$ws.Range("O21").Value = "Región de Ñuble"

# Row 22
$ws.Range("D22").Value = 44517
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 850
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 850

# Row 23
$ws.Range("D23").Value = 44900
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 900
$ws.Range("L23").Value = 1000
$ws.Range("M23").Value = 950
$ws.Range("O23").Value = "Región de Ñuble"
$ws.Range("P23").Value = 950

# Row 24
$ws.Range("D24").Value = 44505
$ws.Range("J24").Value = 440
$ws.Range("K24").Value = 900
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = 950
$ws.Range("P24").Value = 950

# Row 25
$ws.Range("D25").Value = 44530
$ws.Range("J25").Value = 300

# Row 26
$ws.Range("D26").Value = 44880
$ws.Range("J26").Value = 1200
$ws.Range("L26").Value = 1100
$ws.Range("M26").Value = 1050
$ws.Range("O26").Value = "Región de Ñuble"
$ws.Range("P26").Value = 1050

# Row 27
$ws.Range("D27").Value = 44510
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 950
$ws.Range("P27").Value = 950

# Row 28
$ws.Range("D28").Value = 44895
$ws.Range("J28").Value = 1200
$ws.Range("K28").Value = 900
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 950
$ws.Range("O28").Value = "Región de Ñuble"
$ws.Range("P28").Value = 950

# Row 29
$ws.Range("D29").Value = 44860
$ws.Range("L29").Value = 1200
$ws.Range("M29").Value = 1100
$ws.Range("O29").Value = "Provincia de Diguillín"
$ws.Range("P29").Value = 1100

# Row 30
$ws.Range("D30").Value = 44848
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 1300
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1400
$ws.Range("P30").Value = 1400

# Row 31
$ws.Range("D31").Value = 44889
$ws.Range("J31").Value = 600
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = 1000
$ws.Range("P31").Value = 1000

# Row 32
$ws.Range("D32").Value = 44518
$ws.Range("J32").Value = 400
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 900
$ws.Range("M32").Value = 850
$ws.Range("P32").Value = 850

# Row 33
$ws.Range("D33").Value = 44902
$ws.Range("J33").Value = 2000
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = 950
$ws.Range("P33").Value = 950

# Row 35
$ws.Range("D35").Value = 44831
$ws.Range("J35").Value = 200
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 2200
$ws.Range("M35").Value = 2100
$ws.Range("O35").Value = "Provincia de Diguillín"
$ws.Range("P35").Value = 2100

# Row 36
$ws.Range("D36").Value = 44524
$ws.Range("J36").Value = 400

# Row 37
$ws.Range("D37").Value = 44883
$ws.Range("J37").Value = 1000
$ws.Range("O37").Value = "Región de Ñuble"

# Row 38
$ws.Range("D38").Value = 44537
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 800
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = 850
$ws.Range("P38").Value = 850

# Row 39
$ws.Range("D39").Value = 44525
$ws.Range("J39").Value = 360
$ws.Range("K39").Value = 800
$ws.Range("L39").Value = 900
$ws.Range("M39").Value = 850
$ws.Range("P39").Value = 850

# Row 40
$ws.Range("D40").Value = 44890
$ws.Range("J40").Value = 160

# Row 41
$ws.Range("D41").Value = 44827
$ws.Range("J41").Value = 120
$ws.Range("K41").Value = 2200
$ws.Range("L41").Value = 2300
$ws.Range("M41").Value = 2250
$ws.Range("P41").Value = 2250

# Row 42
$ws.Range("D42").Value = 44894

# Row 43
$ws.Range("D43").Value = 44910
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = 950
$ws.Range("O43").Value = "Provincia de Diguillín"
$ws.Range("P43").Value = 950

# Row 44
$ws.Range("D44").Value = 44875
$ws.Range("J44").Value = 1000
$ws.Range("K44").Value = 1000
$ws.Range("L44").Value = 1100
$ws.Range("M44").Value = 1050
$ws.Range("O44").Value = "Provincia de Diguillín"
$ws.Range("P44").Value = 1050

# Row 45
$ws.Range("D45").Value = 44553
$ws.Range("J45").Value = 8000

# Row 46
$ws.Range("D46").Value = 44908
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = 950
$ws.Range("P46").Value = 950

# Row 47
$ws.Range("D47").Value = 44858
$ws.Range("J47").Value = 1000
$ws.Range("K47").Value = 1000
$ws.Range("L47").Value = 1200
$ws.Range("M47").Value = 1100
$ws.Range("P47").Value = 1100
